$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append row 4 (Testmail #3)
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A4").Value = "Wil je dit voor me oppakken?"
$wsLogs.Range("B4").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C4").Value = "Testmail #3: Wil je dit voor me oppakken?"
$wsLogs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$wsLogs.Range("E4").Value = "Beste klant,`nBedankt voor uw bericht. Kunt u meer details geven over wat precies moet worden opgepakt en welke specifieke actie er van ons wordt verwacht? Met meer informatie kan ik u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent"
$wsLogs.Range("F4").Value = "2025-07-27 18:34:55"
$wsLogs.Range("G4").Value = "Ja"
$wsLogs.Range("H4").Value = "Nee"
$wsLogs.Range("I4").Value = "Ja"
$wsLogs.Range("J4").Value = "Nee"

# Writing the multi-line E4 text auto-expands the row height; AutoFit
# brings it back to the (non-custom) default so no ht/customHeight
# override is persisted, matching the un-styled rows 2/3.
$wsLogs.Rows.Item(4).AutoFit()

# Extend the conditional-formatting ranges so row 4 is covered too
# (each block of cfRules shares one sqref; updating the first rule's
# AppliesTo range updates the whole group).
$wsLogs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D4"))
$wsLogs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G4"))
$wsLogs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H4"))
$wsLogs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I4"))
$wsLogs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J4"))

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: append row 3 (new category tally)
# ---------------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$wsDash.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend category/value series refs
#    from Dashboard!$A$2 / $B$2 to $A$2:$A$3 / $B$2:$B$3
# ---------------------------------------------------------------------
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = '=SERIES(''Dashboard''!$B$1,''Dashboard''!$A$2:$A$3,''Dashboard''!$B$2:$B$3,1)'
